$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") ---
# Match the look of the existing header row (bold font, thin border,
# centered horizontal / top vertical alignment) by copying the format
# from the neighboring header cell H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

foreach ($addr in @("I1", "J1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
    $cell.Borders.Weight = 2            # xlThin
}

# --- New data values for columns I ("I0") and J ("IF"), rows 2-32 ---
$values = @(
    @(4, 5),
    @(7, 8),
    @(6, 7),
    @(5, 7),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 2),
    @(1, 7),
    @(1, 7),
    @(1, 7),
    @(1, 7),
    @(1, 5),
    @(1, 3),
    @(1, 4),
    @(1, 7),
    @(1, 7),
    @(1, 4),
    @(1, 5),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 7),
    @(1, 4),
    @(1, 4),
    @(1, 3),
    @(1, 2)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
    $row++
}
